$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("COT_Data_input_sheet")
$ws.Activate()

# --- Data updates (MinEnrolled / MaxEnrolled columns) ---
$ws.Range("D3").Value = 1

$ws.Range("D4").Value = 2
$ws.Range("E4").Value = 9998

$ws.Range("D5").Value = 3
$ws.Range("E5").Value = 9997

$ws.Range("D6").Value = 4
$ws.Range("E6").Value = 9996

$ws.Range("D7").Value = ""
$ws.Range("E7").Value = 99995

# --- View / selection updates ---
# Scroll the sheet so row 1 (B1) is the top-left visible cell, then move the
# selection down to D8, matching the saved workbook view state.
$ws.Range("B1").Select()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("D8").Select()
